$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.684.61"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "3.085.97"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").Value = "3.618.46"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "58.701.60"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "3.083.54"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "0.0₃0922"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "3.126.85"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  +5.81%  "
$ws.Range("D45").Value = "2.282.08"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.957"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.748"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "263.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.53%  "
